$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.985.64"
$ws.Range("E2").Value = "  +2.93%  "
$ws.Range("D3").Value = "1.599.91"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.49"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.06"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").Value = "1.822.36"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "1.594.67"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "25.990.83"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.28"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "0.0₃0721"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "200.87"
$ws.Range("E20").Value = "  +8.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.22"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.00"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("E24").Value = "  +12.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.47"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.123"
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.13"
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.41"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0165"
$ws.Range("E36").Value = "  +11.23%  "
$ws.Range("D37").Value = "1.126.56"
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.791"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.490"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.782"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Value = "1.735.59"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.14"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.06"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.51"
$ws.Range("E46").Value = "  +3.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.31"
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "0.0₇0932"
$ws.Range("E51").Value = "  -16.35%  "
